$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# --- Extend the "Table2" ListObject with a new "Reg Proc" column (T) ---
$lo = $ws.ListObjects.Item("Table2")
$newCol = $lo.ListColumns.Add()

# --- Set the values for the new column / touched cells in the exact order
#     the original author typed them, so new shared-strings come out in the
#     same order as the source workbook. ---
$ws.Range("T2").Value = "Reg Proc"
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"
$ws.Range("T6").Value = "ID Repo- need to know "
$ws.Range("T9").Value = "Under processing`nProcessed"
$ws.Range("T8").Value = "Under processing`nProcessed`n"
$ws.Range("S8").Value = "Reg proc`nArchival policy"
$ws.Range("T10").Value = "E-UIN Generation"
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"

# --- Match formatting of the new column to the look of the neighbouring
#     "Research info" / "Module Dependency" columns it was modelled on. ---
$ws.Range("S2").Copy() | Out-Null
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("S5").Copy() | Out-Null
$ws.Range("T4").PasteSpecial(-4122) | Out-Null
$ws.Range("T5").PasteSpecial(-4122) | Out-Null
$ws.Range("T7").PasteSpecial(-4122) | Out-Null
$ws.Range("T8").PasteSpecial(-4122) | Out-Null
$ws.Range("T9").PasteSpecial(-4122) | Out-Null
$ws.Range("S8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("S6").Copy() | Out-Null
$ws.Range("T6").PasteSpecial(-4122) | Out-Null
$ws.Range("T10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column header / width ---
$newCol.Name = "Reg Proc"
$ws.Columns.Item(20).ColumnWidth = 31.1

# --- Put the cursor where the author left it ---
$ws.Range("T4").Select()
